# Edit script: update column H (Inflow/model temperature) data values for rows 4-111
# on the "Statistics calculator" worksheet. This mirrors a refresh of upstream
# skill-assessment source data; all dependent formulas (J, K, L, M, N, O, P, Q, R,
# H2, B2, B3, B4, etc.) recalculate automatically from these raw inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics calculator")

$newH = @{
    4 = 5.6333679999999999
    5 = 5.8525710000000002
    6 = 5.9628100000000002
    7 = 6.5337459999999998
    8 = 8.5463539999999991
    9 = 12.305540000000001
    10 = 17.797739
    11 = 17.887530999999999
    12 = 14.438943
    13 = 10.489164000000001
    14 = 5.3780950000000001
    15 = 3.7197469999999999
    16 = 4.213533
    17 = 3.3110390000000001
    18 = 4.7018979999999999
    19 = 4.7580549999999997
    20 = 8.2818780000000007
    21 = 13.595072999999999
    22 = 16.21734
    23 = 17.771920999999999
    24 = 17.381886000000002
    25 = 10.232779000000001
    26 = 4.8344560000000003
    27 = 4.0301580000000001
    28 = 3.4191470000000002
    29 = 3.9419179999999998
    30 = 3.101648
    31 = 8.2963290000000001
    32 = 10.669551999999999
    33 = 12.462146000000001
    34 = 17.817267999999999
    35 = 18.962933
    36 = 16.609960999999998
    37 = 10.76773
    38 = 6.8409709999999997
    39 = 2.559825
    40 = 3.1383709999999998
    41 = 3.2248890000000001
    42 = 6.1851079999999996
    43 = 7.8967999999999998
    44 = 11.922974
    45 = 14.898326000000001
    46 = 18.877089000000002
    47 = 19.143405999999999
    48 = 15.331623
    49 = 9.8799829999999993
    50 = 5.785768
    51 = 3.7002160000000002
    52 = 6.1772119999999999
    53 = 3.8808240000000001
    54 = 6.1220439999999998
    55 = 8.2367699999999999
    56 = 12.453976000000001
    57 = 14.089005999999999
    58 = 19.526872999999998
    59 = 19.820018999999998
    60 = 16.607327000000002
    61 = 12.605307
    62 = 6.3983939999999997
    63 = 5.3554399999999998
    64 = 7.7686400000000004
    65 = 8.0858450000000008
    66 = 9.4184230000000007
    67 = 8.124053
    68 = 12.881795
    69 = 18.613572999999999
    70 = 20.134067999999999
    71 = 19.531728999999999
    72 = 14.65448
    73 = 13.600369000000001
    74 = 4.9770370000000002
    75 = 3.8891260000000001
    76 = 3.7530060000000001
    77 = 6.3609549999999997
    78 = 5.7115340000000003
    79 = 10.877504
    80 = 12.746848
    81 = 15.079478999999999
    82 = 16.463688000000001
    83 = 18.133773999999999
    84 = 13.475885999999999
    85 = 8.9190860000000001
    86 = 7.6639989999999996
    87 = 1.808792
    88 = 2.705409
    89 = 4.0119020000000001
    90 = 5.369224
    91 = 6.1563179999999997
    92 = 11.897114
    93 = 15.144710999999999
    94 = 18.398571
    95 = 20.152441
    96 = 15.909276999999999
    97 = 9.7026559999999993
    98 = 5.4410220000000002
    99 = 4.0441940000000001
    100 = 5.1618069999999996
    101 = 3.9569990000000002
    102 = 4.7896400000000003
    103 = 7.6771289999999999
    104 = 13.463158
    105 = 14.901141000000001
    106 = 19.623201000000002
    107 = 18.838439999999999
    108 = 13.972019
    109 = 11.201335
    110 = 6.6556240000000004
    111 = 3.397316
}

foreach ($row in $newH.Keys) {
    $ws.Cells.Item([int]$row, 8).Value = $newH[$row]
}

# Update the active selection left behind by the editor to reflect the edited range
$ws.Activate()
$ws.Range("H3:I111").Select()

$wb.Save()
